# "Texto pronto, agora é ensaiar"
# Fill in the rehearsed slide durations (column F, seconds) for slides
# that had not yet been timed (rows 56-90). Columns G-K recalculate
# automatically from the existing shared formulas.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$durations = [ordered]@{
    56 = 37
    57 = 29
    58 = 56
    59 = 8
    60 = 7
    61 = 13
    62 = 25
    63 = 14
    64 = 42
    65 = 25
    66 = 20
    67 = 13
    68 = 20
    69 = 20
    70 = 20
    71 = 20
    72 = 20
    73 = 10
    74 = 12
    75 = 23
    76 = 20
    77 = 25
    78 = 20
    79 = 27
    80 = 25
    81 = 36
    82 = 0
    83 = 116
    84 = 45
    85 = 17
    86 = 11
    87 = 8
    88 = 10
    89 = 10
    90 = 10
}

foreach ($row in $durations.Keys) {
    $ws.Cells.Item($row, 6).Value = $durations[$row]
}

# Leave the view where it ended up after entering the last value.
$ws.Range("F91").Select()
